$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.130774
$ws.Range("B2").Value = 0.6592673999999999
$ws.Range("C2").Value = 0.130774
$ws.Range("D2").Value = 0.7927536186445832
$ws.Range("E2").Value = 0.7927536186445832
$ws.Range("F2").Value = 1

$ws.Range("A3").Value = 0.153969
$ws.Range("B3").Value = 0.6592673999999999
$ws.Range("C3").Value = 0.02319499999999999
$ws.Range("D3").Value = 0.7927536186445832
$ws.Range("E3").Value = 0.7927536186445832
$ws.Range("F3").Value = 0.1506472082042489

$ws.Range("A4").Value = 0.182164
$ws.Range("B4").Value = 0.6589341
$ws.Range("C4").Value = 0.028195
$ws.Range("D4").Value = 0.792352833195319
$ws.Range("E4").Value = 0.7901641992234182
$ws.Range("F4").Value = 0.1547781120309172

$ws.Range("A5").Value = 0.210358
$ws.Range("B5").Value = 0.6576009
$ws.Range("C5").Value = 0.028194
$ws.Range("D5").Value = 0.7907496913982623
$ws.Range("E5").Value = 0.7901641992234182
$ws.Range("F5").Value = 0.1340286559104004

$ws.Range("A6").Value = 0.244951
$ws.Range("B6").Value = 0.6556010999999999
$ws.Range("C6").Value = 0.03459300000000001
$ws.Range("D6").Value = 0.7883449787026769
$ws.Range("E6").Value = 0.7772824056300296
$ws.Range("F6").Value = 0.141224163199987

$ws.Range("A7").Value = 0.278144
$ws.Range("B7").Value = 0.6526013999999999
$ws.Range("C7").Value = 0.033193
$ws.Range("D7").Value = 0.7847379096592991
$ws.Range("E7").Value = 0.7772824056300296
$ws.Range("F7").Value = 0.1193374654855039

$ws.Range("A8").Value = 0.314137
$ws.Range("B8").Value = 0.6496016999999999
$ws.Range("C8").Value = 0.035993
$ws.Range("D8").Value = 0.7811308406159214
$ws.Range("E8").Value = 0.7772824056300296
$ws.Range("F8").Value = 0.1145773977595762

$ws.Range("A9").Value = 0.342332
$ws.Range("B9").Value = 0.646602
$ws.Range("C9").Value = 0.02819500000000003
$ws.Range("D9").Value = 0.7775237715725437
$ws.Range("E9").Value = 0.7772824056300296
$ws.Range("F9").Value = 0.08236156713365979

$ws.Range("A10").Value = 0.378324
$ws.Range("B10").Value = 0.643269
$ws.Range("C10").Value = 0.03599199999999997
$ws.Range("D10").Value = 0.7735159170799017
$ws.Range("E10").Value = 0.737691579495373
$ws.Range("F10").Value = 0.09513538659984555

$ws.Range("A11").Value = 0.408918
$ws.Range("B11").Value = 0.6399359999999999
$ws.Range("C11").Value = 0.03059400000000001
$ws.Range("D11").Value = 0.7695080625872597
$ws.Range("E11").Value = 0.737691579495373
$ws.Range("F11").Value = 0.07481695596672196

$ws.Range("A12").Value = 0.441112
$ws.Range("B12").Value = 0.636603
$ws.Range("C12").Value = 0.032194
$ws.Range("D12").Value = 0.7655002080946178
$ws.Range("E12").Value = 0.737691579495373
$ws.Range("F12").Value = 0.07298373202270625

$ws.Range("A13").Value = 0.471906
$ws.Range("B13").Value = 0.6336033
$ws.Range("C13").Value = 0.03079399999999999
$ws.Range("D13").Value = 0.76189313905124
$ws.Range("E13").Value = 0.737691579495373
$ws.Range("F13").Value = 0.06525452102749274

$ws.Range("A14").Value = 0.5001
$ws.Range("B14").Value = 0.6306036
$ws.Range("C14").Value = 0.028194
$ws.Range("D14").Value = 0.7582860700078623
$ws.Range("E14").Value = 0.737691579495373
$ws.Range("F14").Value = 0.05637672465506898

$ws.Range("A15").Value = 0.530694
$ws.Range("B15").Value = 0.6276039
$ws.Range("C15").Value = 0.03059400000000001
$ws.Range("D15").Value = 0.7546790009644845
$ws.Range("E15").Value = 0.737691579495373
$ws.Range("F15").Value = 0.05764904069011523

$ws.Range("A16").Value = 0.565487
$ws.Range("B16").Value = 0.6239376
$ws.Range("C16").Value = 0.03479299999999996
$ws.Range("D16").Value = 0.7502703610225783
$ws.Range("E16").Value = 0.737691579495373
$ws.Range("F16").Value = 0.0615274975375207

$ws.Range("A17").Value = 0.606479
$ws.Range("B17").Value = 0.619938
$ws.Range("C17").Value = 0.04099200000000003
$ws.Range("D17").Value = 0.7454609356314079
$ws.Range("E17").Value = 0.737691579495373
$ws.Range("F17").Value = 0.06759013914743961

$ws.Range("A18").Value = 0.64867
$ws.Range("B18").Value = 0.6156050999999999
$ws.Range("C18").Value = 0.04219099999999998
$ws.Range("D18").Value = 0.7402507247909732
$ws.Range("E18").Value = 0.737691579495373
$ws.Range("F18").Value = 0.06504231735705364

$ws.Range("A19").Value = 0.689662
$ws.Range("B19").Value = 0.6119387999999999
$ws.Range("C19").Value = 0.04099200000000003
$ws.Range("D19").Value = 0.735842084849067
$ws.Range("E19").Value = 0.7065751134346605
$ws.Range("F19").Value = 0.05943781156566554

$ws.Range("A20").Value = 0.737253
$ws.Range("B20").Value = 0.6072725999999999
$ws.Range("C20").Value = 0.04759100000000005
$ws.Range("D20").Value = 0.7302310885593684
$ws.Range("E20").Value = 0.7065751134346605
$ws.Range("F20").Value = 0.06455178887030646

$ws.Range("A21").Value = 0.787243
$ws.Range("B21").Value = 0.6026063999999999
$ws.Range("C21").Value = 0.04998999999999998
$ws.Range("D21").Value = 0.7246200922696695
$ws.Range("E21").Value = 0.7065751134346605
$ws.Range("F21").Value = 0.06350008828277924

$ws.Range("A22").Value = 0.835833
$ws.Range("B22").Value = 0.5986068
$ws.Range("C22").Value = 0.04859000000000002
$ws.Range("D22").Value = 0.7198106668784993
$ws.Range("E22").Value = 0.7065751134346605
$ws.Range("F22").Value = 0.05813362238629011

$ws.Range("A23").Value = 0.8808240000000001
$ws.Range("B23").Value = 0.5946071999999999
$ws.Range("C23").Value = 0.044991
$ws.Range("D23").Value = 0.7150012414873287
$ws.Range("E23").Value = 0.7065751134346605
$ws.Range("F23").Value = 0.0510783084929566

$ws.Range("A24").Value = 0.925615
$ws.Range("B24").Value = 0.5909409
$ws.Range("C24").Value = 0.04479099999999991
$ws.Range("D24").Value = 0.7105926015454227
$ws.Range("E24").Value = 0.7065751134346605
$ws.Range("F24").Value = 0.04839052953981938

$ws.Range("A25").Value = 0.960208
$ws.Range("B25").Value = 0.5876079
$ws.Range("C25").Value = 0.03459299999999998
$ws.Range("D25").Value = 0.7065847470527806
$ws.Range("E25").Value = 0.7065751134346605
$ws.Range("F25").Value = 0.0360265692433306

$ws.Range("A26").Value = 1
$ws.Range("B26").Value = 0.5849415
$ws.Range("C26").Value = 0.03979200000000005
$ws.Range("D26").Value = 0.7033784634586672
$ws.Range("E26").Value = 0.6262411273069634
$ws.Range("F26").Value = 0.03979200000000005

